$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-shape the schema -------------------------------------------------
# Old layout: ticket_id, created_at, source, category, sub_category,
#             severity, raw_incident, AI_Incident_summary, assigned_engineer,
#             sla_hours, status, Engineer_Updates, closing_time
# New layout: ticket_id, created_at, source, raw_incident,
#             AI_Incident_summary, assigned_engineer, status,
#             Engineer_Updates, closing_time
# Drop sla_hours (J), then category/sub_category/severity (D:F). Deleting
# entire columns shifts everything left and keeps header text/styles intact.
$ws.Columns("J:J").Delete()
$ws.Columns("D:F").Delete()

# The column delete shifts the (already blank) F/H/I cells of rows 2-3 left
# as empty-string text cells instead of truly-blank ones; clear them back to
# plain blanks so they read the same way as the rest of the blank cells in
# the sheet (e.g. the new rows added below).
$ws.Range("F2:F3").ClearContents()
$ws.Range("H2:I3").ClearContents()

# --- Add two more incident rows so rows 2-5 exist ------------------------
$ws.Rows("4:5").Insert()

# --- Row 2 ----------------------------------------------------------------
$ws.Range("A2").Value = "TCK-EFABAF75"
$ws.Range("B2").Value = "2026-01-02 09:55:45"
$ws.Range("C2").Value = "Email"
$ws.Range("D2").Value = "Incident Title:`nIncident Description:`nGross sales margin discrepancy in Power BI report SEF_ALL (production). Observed today. Business impact: Director needs accurate sales figures for closing period. Urgent."
$ws.Range("E2").Value = "Summarized Mail in Description"
$ws.Range("G2").Value = "OPEN"

# --- Row 3 ----------------------------------------------------------------
$ws.Range("A3").Value = "TCK-677FB46A"
$ws.Range("B3").Value = "2026-01-02 10:13:38"
$ws.Range("C3").Value = "Email"
$ws.Range("D3").Value = "Incident Title:`nIncident Description:`nGross sales margin discrepancy in Power BI report SEF_ALL (production). Observed today. Business impact: Director needs accurate sales figures for closing period. Urgent."
$ws.Range("E3").Value = "Summarized Mail in Description"
$ws.Range("G3").Value = "OPEN"

# --- Row 4 (new) ------------------------------------------------------------
$ws.Range("A4").Value = "TCK-48A837EA"
$ws.Range("B4").Value = "2026-01-02 10:13:41"
$ws.Range("C4").Value = "Email"
$ws.Range("D4").Value = "Incident Title:`nIncident Description:`nGross sales margin discrepancy in Power BI report SEF_ALL (production). Observed today. Business impact: Director needs accurate sales figures for closing period. Urgent."
$ws.Range("E4").Value = "Summarized Mail in Description"
$ws.Range("G4").Value = "OPEN"

# --- Row 5 (new) ------------------------------------------------------------
$ws.Range("A5").Value = "TCK-2F7D24EA"
$ws.Range("B5").Value = "2026-01-02 10:41:42"
$ws.Range("C5").Value = "Email"
$ws.Range("D5").Value = "Incident Title:`nIncident Description:`nIncident: Gross sales margin discrepancy in Consolidated Sales Report. System: Azure Synapse. Environment: Production. Noticed this morning. Business impact: None reported. Urgency: Director needs figures for closing period. Request immediate investigation."
$ws.Range("E5").Value = "Summarized Mail in Description"
$ws.Range("G5").Value = "OPEN"

# Entering the multi-line incident text auto-expands row height; put it back
# to the sheet's default (un-set any "custom height" marker) to match the
# original, non-autofitted layout.
$ws.Rows("2:5").AutoFit()
